$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4913134422460475
$ws.Range("C2").Value = 0.9902169060012099
$ws.Range("D2").Value = 0.5666452762251486
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(n_estimators=100))])"
$ws.Range("G2").Value = 0.1209129460333012
$ws.Range("H2").Value = 0.99
